$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI ligand-receptor metrics (Vegfb -> Nrp1) following Dr Hou advice.
# Each row below corresponds to one data row (sheet row number) in column A,
# followed by the new values for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T (in that order).
$updates = @(
    ,(2, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 2, 107.663086, 215.326172, 0.2751823527645522, 0.2057131686791961, 140.330596374921, 561.322385499684, 0.01381272833466232, 0.007660980575213401)
    ,(3, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 3, 27.41172733333334, 82.23518200000001, 0.07006323059444025, 0.07856388152449204, 35.72908958185901, 214.374537491154, 0.003516811164405699, 0.002925803798254208)
    ,(4, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 3, 82.303927, 246.911781, 0.2103654011331419, 0.2358886724356653, 107.2768725940845, 643.661235564507, 0.01055925319219327, 0.008784748925143875)
    ,(5, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 3, 71.272429, 213.817287, 0.1821693528222338, 0.2042716462128862, 92.8981588606815, 557.388953164089, 0.009143957656280708, 0.007607296721700085)
    ,(6, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 3, 83.25665266666667, 249.769958, 0.2128005281598046, 0.2386192492246123, 108.518677617071, 651.1120657024261, 0.01068148395205768, 0.00888643855383203)
    ,(7, 2, 1.3034235, 2.606847, 0.05019481880249994, 0.03724108001642075, 2, 19.3348755, 38.66975100000001, 0.04941913452582716, 0.03694338192314826, 25.20153109627426, 100.806124385097, 0.002480584502900264, 0.001375811442277156)
    ,(8, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 2, 107.663086, 215.326172, 0.2751823527645522, 0.2057131686791961, 439.9331532270673, 2639.598919362404, 0.04330258181688194, 0.03602549367345879)
    ,(9, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 3, 27.41172733333334, 82.23518200000001, 0.07006323059444025, 0.07856388152449204, 112.0098641901638, 1008.088777711474, 0.01102512114127702, 0.01375849020748269)
    ,(10, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 3, 82.303927, 246.911781, 0.2103654011331419, 0.2358886724356653, 336.3104985498964, 3026.794486949067, 0.03310301297482945, 0.04130997510287761)
    ,(11, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 3, 71.272429, 213.817287, 0.1821693528222338, 0.2042716462128862, 291.2335656821343, 2621.102091139209, 0.0286660944129022, 0.03577304722667258)
    ,(12, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 3, 83.25665266666667, 249.769958, 0.2128005281598046, 0.2386192492246123, 340.2035283920562, 3061.831755528506, 0.0334862035618973, 0.04178816703131224)
    ,(13, 3, 4.086202333333333, 12.258607, 0.1573595885849988, 0.1751248785129528, 2, 19.3348755, 38.66975100000001, 0.04941913452582716, 0.03694338192314826, 79.00621338280952, 474.0372802968571, 0.00777657467721087, 0.006469705271148958)
    ,(14, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 2, 107.663086, 215.326172, 0.2751823527645522, 0.2057131686791961, 584.0111248048472, 3504.066748829084, 0.05748416396519761, 0.04782383170614865)
    ,(15, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 3, 27.41172733333334, 82.23518200000001, 0.07006323059444025, 0.07856388152449204, 148.6930600594504, 1338.237540535054, 0.01463584490415427, 0.0182643914939077)
    ,(16, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 3, 82.303927, 246.911781, 0.2103654011331419, 0.2358886724356653, 446.4520827791063, 4018.068745011957, 0.04394423948285911, 0.05483897916881852)
    ,(17, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 3, 71.272429, 213.817287, 0.1821693528222338, 0.2042716462128862, 386.6124683428043, 3479.512215085239, 0.03805423146457, 0.04748870912614044)
    ,(18, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 3, 83.25665266666667, 249.769958, 0.2128005281598046, 0.2386192492246123, 451.6200786901695, 4064.580708211526, 0.04445292486864229, 0.05547377880587511)
    ,(19, 3, 5.424432333333333, 16.273297, 0.2088948051635471, 0.232478222046779, 2, 19.3348755, 38.66975100000001, 0.04941913452582716, 0.03694338192314826, 104.8807238231745, 629.2843429390471, 0.01032340047812379, 0.008588531745888623)
    ,(20, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 2, 107.663086, 215.326172, 0.2751823527645522, 0.2057131686791961, 715.6208856068347, 4293.725313641008, 0.07043850121671823, 0.05860116587693855)
    ,(21, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 3, 27.41172733333334, 82.23518200000001, 0.07006323059444025, 0.07856388152449204, 182.2017677469609, 1639.815909722648, 0.01793410407278636, 0.02238036136778688)
    ,(22, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 3, 82.303927, 246.911781, 0.2103654011331419, 0.2358886724356653, 547.0622412649427, 4923.560171384484, 0.05384728858812563, 0.06719721110052207)
    ,(23, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 3, 71.272429, 213.817287, 0.1821693528222338, 0.2042716462128862, 473.7374773033187, 4263.637295729868, 0.04662993848081749, 0.05819052178591638)
    ,(24, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 3, 83.25665266666667, 249.769958, 0.2128005281598046, 0.2386192492246123, 553.3948703084792, 4980.553832776312, 0.05447060871133572, 0.0679750659377996)
    ,(25, 3, 6.646854666666667, 19.940564, 0.2559702702919538, 0.2848683254124845, 2, 19.3348755, 38.66975100000001, 0.04941913452582716, 0.03694338192314826, 128.516107446594, 771.0966446795642, 0.01264982922217041, 0.0105239993435211)
    ,(26, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 2, 107.663086, 215.326172, 0.2751823527645522, 0.2057131686791961, 205.3306593203693, 1231.983955922216, 0.02021067884862366, 0.01681423260341355)
    ,(27, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 3, 27.41172733333334, 82.23518200000001, 0.07006323059444025, 0.07856388152449204, 52.27853162655511, 470.506784638996, 0.005145771298252022, 0.006421520735213031)
    ,(28, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 3, 82.303927, 246.911781, 0.2103654011331419, 0.2358886724356653, 156.9667025480353, 1412.700322932318, 0.0154502188110934, 0.01928066653345376)
    ,(29, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 3, 71.272429, 213.817287, 0.1821693528222338, 0.2042716462128862, 135.9278781766873, 1223.350903590186, 0.01337936916725879, 0.01669640789531536)
    ,(30, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 3, 83.25665266666667, 249.769958, 0.2128005281598046, 0.2386192492246123, 158.7837021953249, 1429.053319757924, 0.01562906592754927, 0.01950385376821186)
    ,(31, 3, 1.907159333333333, 5.721477999999999, 0.07344467639578636, 0.08173629676394162, 2, 19.3348755, 38.66975100000001, 0.04941913452582716, 0.03694338192314826, 36.874688268663, 221.248129611978, 0.003629572343009209, 0.003019615228334083)
    ,(32, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 2, 107.663086, 215.326172, 0.2751823527645522, 0.2057131686791961, 710.492336561377, 2841.969346245508, 0.06993369858246846, 0.03878746424402314)
    ,(33, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 3, 27.41172733333334, 82.23518200000001, 0.07006323059444025, 0.07856388152449204, 180.8960055468164, 1085.376033280898, 0.01780557801356489, 0.01481331392184753)
    ,(34, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 3, 82.303927, 246.911781, 0.2103654011331419, 0.2358886724356653, 543.1416799849766, 3258.850079909859, 0.05346138808404107, 0.04447709160484947)
    ,(35, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 3, 71.272429, 213.817287, 0.1821693528222338, 0.2042716462128862, 470.3424032691655, 2822.054419614993, 0.04629576164040463, 0.03851566345714135)
    ,(36, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 3, 83.25665266666667, 249.769958, 0.2128005281598046, 0.2386192492246123, 549.4289257825936, 3296.573554695562, 0.05408024113832235, 0.04499194512758144)
    ,(37, 2, 6.5992195, 13.198439, 0.254135840761214, 0.1885511972474212, 2, 19.3348755, 38.66975100000001, 0.04941913452582716, 0.03694338192314826, 127.5950874296723, 510.3803497186891, 0.01255917330241263, 0.006965718891978344)
)

$cols = @(5, 7, 8, 9, 10, 11, 13, 14, 15, 16, 17, 18, 19, 20)  # E,G,H,I,J,K,M,N,O,P,Q,R,S,T

foreach ($update in $updates) {
    $row = $update[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $update[$i + 1]
    }
}

